$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New patient row (row 2): MARTIN / ANDRE / Marie, born 05/05/1940, sexe F
$ws.Range("A2").Value = "MARTIN"
$ws.Range("B2").Value = "ANDRE"
$ws.Range("C2").Value = "Marie"
$ws.Range("D2").Value = 14736
$ws.Range("U2").Value = "F"

# New patient row (row 3): BERNARD / Jean
$ws.Range("A3").Value = "BERNARD"
$ws.Range("C3").Value = "Jean"

# Fix birth-date storage: replace placeholder text dates with the real
# numeric date (05/05/1940) on the other rows using the same value
$ws.Range("D5").Value = 14736
$ws.Range("D8").Value = 14736
$ws.Range("D11").Value = 14736
$ws.Range("D14").Value = 14736
$ws.Range("D17").Value = 14736

# Update the active selection to reflect where editing finished
$ws.Range("T5").Select()
